$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4314.0586
$ws.Range("I62").Value = 2083.6875
$ws.Range("J62").Value = 40000
$ws.Range("K62").Value = 2083.6875
$ws.Range("L62").Value = 40000
$ws.Range("M62").Value = -1459.6875
$ws.Range("N62").Value = -41248

$ws.Range("H65").Value = 4314.0586
$ws.Range("I65").Value = 2083.6875
$ws.Range("J65").Value = 40000
$ws.Range("K65").Value = 10418.4375
$ws.Range("L65").Value = 200000
$ws.Range("M65").Value = -7298.4375
$ws.Range("N65").Value = -206240

$ws.Range("H133").Value = 74780
$ws.Range("J133").Value = 74780
$ws.Range("L133").Value = 74780
$ws.Range("N133").Value = -84900

$ws.Range("H137").Value = 1756.8837
$ws.Range("I137").Value = 1488.0526
$ws.Range("K137").Value = 4464.1578
$ws.Range("M137").Value = -1914.1578

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 41
$ws.Range("I4").Value = 41
$ws.Range("K4").Value = 41
$ws.Range("M4").Value = 75

$ws.Range("H5").Value = 437.18182
$ws.Range("I5").Value = 400.77777
$ws.Range("J5").Value = 601
$ws.Range("K5").Value = 400.77777
$ws.Range("L5").Value = 601
$ws.Range("M5").Value = -288.77777
$ws.Range("N5").Value = -825

$ws.Range("H61").Value = 9053.735000000001
$ws.Range("I61").Value = 5115.5415
$ws.Range("J61").Value = 18505.4
$ws.Range("K61").Value = 5115.5415
$ws.Range("L61").Value = 18505.4
$ws.Range("M61").Value = -4903.5415
$ws.Range("N61").Value = -18929.4

$ws.Range("H74").Value = 4766.6113
$ws.Range("I74").Value = 2064.9678
$ws.Range("J74").Value = 21516.8
$ws.Range("K74").Value = 2064.9678
$ws.Range("L74").Value = 21516.8
$ws.Range("M74").Value = -1190.9678
$ws.Range("N74").Value = -23264.8

$ws.Range("H77").Value = 4766.6113
$ws.Range("I77").Value = 2064.9678
$ws.Range("J77").Value = 21516.8
$ws.Range("K77").Value = 10324.839
$ws.Range("L77").Value = 107584
$ws.Range("M77").Value = -5956.839
$ws.Range("N77").Value = -116320

$ws.Range("H132").Value = 2257.7932
$ws.Range("I132").Value = 1445.4736
$ws.Range("K132").Value = 4336.4208
$ws.Range("M132").Value = -1806.4208

$ws.Range("H136").Value = 9053.735000000001
$ws.Range("I136").Value = 5115.5415
$ws.Range("J136").Value = 18505.4
$ws.Range("K136").Value = 15346.6245
$ws.Range("L136").Value = 55516.2
$ws.Range("M136").Value = -12796.6245
$ws.Range("N136").Value = -60616.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 437.18182
$ws.Range("I4").Value = 400.77777
$ws.Range("J4").Value = 601
$ws.Range("K4").Value = 400.77777
$ws.Range("L4").Value = 601
$ws.Range("M4").Value = -285.77777
$ws.Range("N4").Value = -831

$ws.Range("H99").Value = 1298.1818
$ws.Range("I99").Value = 1114.2858
$ws.Range("J99").Value = 1620
$ws.Range("K99").Value = 1114.2858
$ws.Range("L99").Value = 1620
$ws.Range("M99").Value = 383.7141999999999
$ws.Range("N99").Value = -4616

$ws.Range("H134").Value = 35656.332
$ws.Range("I134").Value = 2561.647
$ws.Range("K134").Value = 7684.941
$ws.Range("M134").Value = -5149.941

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6349.222
$ws.Range("I31").Value = 8579.071
$ws.Range("K31").Value = 8579.071
$ws.Range("M31").Value = -8284.071

$ws.Range("H34").Value = 6349.222
$ws.Range("I34").Value = 8579.071
$ws.Range("K34").Value = 8579.071
$ws.Range("M34").Value = -8377.071

$ws.Range("H58").Value = 3136797.5
$ws.Range("I58").Value = 4330393
$ws.Range("J58").Value = 3609.25
$ws.Range("K58").Value = 4330393
$ws.Range("L58").Value = 3609.25
$ws.Range("M58").Value = -4330190
$ws.Range("N58").Value = -4015.25

$ws.Range("H132").Value = 2667.5
$ws.Range("I132").Value = 2260.0667
$ws.Range("K132").Value = 6780.2001
$ws.Range("M132").Value = -4250.2001

$ws.Range("H134").Value = 2497.5881
$ws.Range("I134").Value = 2607.5715
$ws.Range("J134").Value = 2319.923
$ws.Range("K134").Value = 7822.7145
$ws.Range("L134").Value = 6959.768999999999
$ws.Range("M134").Value = -5287.7145
$ws.Range("N134").Value = -12029.769

$ws.Range("H136").Value = 3136797.5
$ws.Range("I136").Value = 4330393
$ws.Range("J136").Value = 3609.25
$ws.Range("K136").Value = 12991179
$ws.Range("L136").Value = 10827.75
$ws.Range("M136").Value = -12988629
$ws.Range("N136").Value = -15927.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 28346.4
$ws.Range("I131").Value = 1746.3636
$ws.Range("J131").Value = 40538.082
$ws.Range("K131").Value = 5239.0908
$ws.Range("L131").Value = 121614.246
$ws.Range("M131").Value = -199.0907999999999
$ws.Range("N131").Value = -131694.246

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3076.92
$ws.Range("I102").Value = 2543.325
$ws.Range("J102").Value = 5211.3
$ws.Range("K102").Value = 2543.325
$ws.Range("L102").Value = 5211.3
$ws.Range("M102").Value = -921.3249999999998
$ws.Range("N102").Value = -8455.299999999999

$ws.Range("H122").Value = 4777.8
$ws.Range("I122").Value = 6874.1113
$ws.Range("J122").Value = 1633.3334
$ws.Range("K122").Value = 20622.3339
$ws.Range("L122").Value = 4900.0002
$ws.Range("M122").Value = -18172.3339
$ws.Range("N122").Value = -9800.0002

$ws.Range("H132").Value = 5080.3423
$ws.Range("I132").Value = 2062.2173
$ws.Range("J132").Value = 9708.134
$ws.Range("K132").Value = 6186.651899999999
$ws.Range("L132").Value = 29124.402
$ws.Range("M132").Value = -3656.651899999999
$ws.Range("N132").Value = -34184.402

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3614.4243
$ws.Range("I132").Value = 3667.0476
$ws.Range("J132").Value = 3522.3333
$ws.Range("K132").Value = 11001.1428
$ws.Range("L132").Value = 10566.9999
$ws.Range("M132").Value = -8471.1428
$ws.Range("N132").Value = -15626.9999

$ws.Range("H136").Value = 5921.091
$ws.Range("I136").Value = 3833.5625
$ws.Range("K136").Value = 11500.6875
$ws.Range("M136").Value = -8950.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9950.916999999999
$ws.Range("I122").Value = 1050.75
$ws.Range("J122").Value = 27751.25
$ws.Range("K122").Value = 3152.25
$ws.Range("L122").Value = 83253.75
$ws.Range("M122").Value = -702.25
$ws.Range("N122").Value = -88153.75

$ws.Range("H132").Value = 2042.3077
$ws.Range("I132").Value = 1341.48
$ws.Range("J132").Value = 3293.7856
$ws.Range("K132").Value = 4024.44
$ws.Range("L132").Value = 9881.356800000001
$ws.Range("M132").Value = -1494.44
$ws.Range("N132").Value = -14941.3568

$ws.Range("H135").Value = 400018530
$ws.Range("J135").Value = 400018530
$ws.Range("L135").Value = 400018530
$ws.Range("N135").Value = -400028670

$ws.Range("H136").Value = 6111.528
$ws.Range("I136").Value = 2620.2
$ws.Range("J136").Value = 10475.6875
$ws.Range("K136").Value = 7860.599999999999
$ws.Range("L136").Value = 31427.0625
$ws.Range("M136").Value = -5310.599999999999
$ws.Range("N136").Value = -36527.0625
